$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).RowHeight = 60
$ws.Rows.Item(39).RowHeight = 60
$ws.Cells.Item(38, 2).Value = "Person Probationer Indicator"
$ws.Cells.Item(38, 5).Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:PersonProbationerIndicator"
$ws.Cells.Item(39, 1).Value = "x-ext"
$ws.Cells.Item(39, 2).Value = "Person Incarcerated Indicator"
$ws.Cells.Item(39, 5).Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/cq-res-ext:PersonIncarceratedIndicator"

for ($r = 36; $r -le 42; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    Write-Host "Row $r - A:[$a] B:[$b] C:[$c] E:[$e]"
}
